$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 50 (shifts existing row 50+ down by one)
$ws.Rows.Item(50).Insert()

$ws.Range("A49:E49").Copy()
$ws.Range("A50:E50").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("B50").Value = "Tạo phíếu mua hàng "

[void]$ws.Range("B50").Select()

Write-Host "done"
